$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.204214811325073
$ws.Range("B1").Value = 2.729588508605957
$ws.Range("C1").Value = 1.714097380638123
$ws.Range("D1").Value = 1.45933735370636
$ws.Range("E1").Value = 1.371708750724792
